# "moving range of cells" ---------------------------------------------------
# 1) The picture on the sheet moves from anchor column B (idx 1) to column
#    H (idx 7).
# 2) A4:B4 (value + formula) moves two columns right, to C4:D4.
# 3) A7:A8 moves down one row, to A8:A9 (this empties/removes row 7 and
#    overwrites what used to be in A9).
# 4) Rows 12:14 end up present (blank) afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Move the picture anchor from col B to col H ------------------------
$shp = $ws.Shapes.Item(1)
$shp.Left = $ws.Cells.Item(6, 8).Left
$shp.Top = $ws.Cells.Item(6, 8).Top

# --- 2. Move A4:B4 -> C4:D4 --------------------------------------------------
# Done manually (formula text copied, then source cleared) because Cut
# collapses a formula cell down to a plain value in this engine.
$ws.Range("D4").Formula = $ws.Range("B4").Formula
$ws.Range("C4").Value = $ws.Range("A4").Value2
$ws.Range("A4").ClearContents()
$ws.Range("B4").ClearContents()

# --- 3. Move A7:A8 -> A8:A9 --------------------------------------------------
# Two single-cell moves, lower one first, so the source/destination overlap
# doesn't clobber data before it is relocated.
$ws.Range("A8").Cut($ws.Range("A9"))
$ws.Range("A7").Cut($ws.Range("A8"))

# --- 4. Leave blank row stubs at rows 12:14 ---------------------------------
$ws.Rows.Item(12).OutlineLevel = 0
$ws.Rows.Item(13).OutlineLevel = 0
$ws.Rows.Item(14).OutlineLevel = 0
